$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 340 (pushes existing rows 340-438 down to 341-439,
# carrying their formatting/styles with them).
$ws.Rows("340:340").Insert()

# Populate the newly inserted row 340 with the new data record.
$ws.Range("A340").Value = 3
$ws.Range("B340").Value = "Femacal de La Calera"
$ws.Range("C340").Value = "Coquimbo"
$ws.Range("D340").Value = 44841
$ws.Range("E340").Value = 5
$ws.Range("F340").Value = 100112040
$ws.Range("G340").Value = "Cilantro"
$ws.Range("H340").Value = "Sin especificar"
$ws.Range("I340").Value = "Primera"
$ws.Range("J340").Value = 120
$ws.Range("K340").Value = 3000
$ws.Range("L340").Value = 3000
$ws.Range("M340").Value = 3000
$ws.Range("N340").Value = "$/docena de atados (3 kilos)"
$ws.Range("O340").Value = "Provincia de Quillota"
$ws.Range("P340").Value = 1000
$ws.Range("Q340").Value = 3
$ws.Range("R340").Value = "Hortaliza"
